$d = $word.ActiveDocument

# 1. Fix page number "p.2" -> "p. 2"
$d.Content.Find.Execute("p.2", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "p. 2", 2)
